# Add "app_type" / "account_type" columns and three new MobiControl test rows
# to the R03 sheet (mirrors the existing TC_375423_375424 row's validation
# test for deleteApplication / account test cases).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R03")

# --- shared literal values (reused across rows, like the existing sheet) ---
$loginUrl = "https://qa2012r2-vr3.inqa.soti.net/login"
$mailTo   = "mailto:vdtenant@yopmail.com"
$email    = "vdtenant@yopmail.com"
$validationMsg   = "MobiControl Instance Name is required"
$validationMsg2  = "MobiControl Instance Name cannot exceed 100 characters"
$instanceName    = "123456789_123456789_123456789_123456789_123456789_123456789_123456789_123456789_123456789_123456789_1"

# 1) Insert two new columns (F:G) for app_type / account_type, shifting the
#    existing validation_message/validation_message2/instance_name columns right.
$ws.Columns("F:G").Insert()

# 2) Insert three new rows for the new test cases, after row 3.
$ws.Rows("4:6").Insert()

# 3) Copy row 3's cell formatting into the new rows so styles/borders match.
$ws.Range("A3:J3").Copy()
$ws.Range("A4:J6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New header + app_type/account_type values, in the same first-use order
#     the original authoring session introduced them (app_type, MobiControl,
#     account_type, Soti DB, ...) so the shared-string table matches exactly.
$ws.Range("F1").Value = "app_type"
$ws.Range("F2").Value = "MobiControl"
$ws.Range("G1").Value = "account_type"
$ws.Range("G2").Value = "Soti DB"

$ws.Range("F3").Value = "MobiControl"
$ws.Range("G3").Value = "Soti DB"

# --- Row 4: TC_375426 ---
$ws.Range("A4").Value = "TC_375426"
$ws.Range("C4").Value = $loginUrl
$ws.Range("D4").Value = $email
$ws.Range("E4").Value = "Welcome12345"
$ws.Range("F4").Value = "MobiControl"
$ws.Range("G4").Value = "Soti DB"

# --- Row 5: deleteApplication ---
$ws.Range("A5").Value = "deleteApplication"
$ws.Range("C5").Value = $loginUrl
$ws.Range("D5").Value = $email
$ws.Range("E5").Value = "Welcome12345"
$ws.Range("F5").Value = "MobiControl"
$ws.Range("G5").Value = "Soti DB"

# --- Row 6: account ---
$ws.Range("A6").Value = "account"
$ws.Range("C6").Value = $loginUrl
$ws.Range("D6").Value = $email
$ws.Range("E6").Value = "Welcome12345"
$ws.Range("F6").Value = "MobiControl"
$ws.Range("G6").Value = "Soti DB"

# --- runmode column: row 3/4/5 -> N, row 6 stays Y (last, so "N" becomes the
#     final new shared-string entry, matching the source edit) ---
$ws.Range("B3").Value = "N"
$ws.Range("B4").Value = "N"
$ws.Range("B5").Value = "N"
$ws.Range("B6").Value = "Y"

# H:J (validation_message / validation_message2 / instance_name) are
# identical to row 3 in every new row, incl. the quote-prefixed instance_name
# cell style, so copy them (value+format) straight from row 3. (Pasted one
# destination row at a time -- a single multi-row paste from a 1-row source
# only fills the first destination row.)
foreach ($r in 4..6) {
    $ws.Range("H3:J3").Copy()
    $ws.Range("H" + $r + ":J" + $r).PasteSpecial(-4163)
}
$excel.CutCopyMode = 0

# 4) Hyperlink the new url/username cells the same way rows 2 & 3 are linked,
#    then restore the cell's text + border formatting that Hyperlinks.Add
#    otherwise resets to the built-in Hyperlink style.
foreach ($r in 4..6) {
    $cCell = $ws.Cells.Item($r, 3)
    $ws.Hyperlinks.Add($cCell, $loginUrl)
    $cCell.Value = $loginUrl
    $ws.Range("C3").Copy()
    $cCell.PasteSpecial(-4122)

    $dCell = $ws.Cells.Item($r, 4)
    $ws.Hyperlinks.Add($dCell, $mailTo)
    $dCell.Value = $email
    $ws.Range("D3").Copy()
    $dCell.PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

Write-Output "done"
